# Time cost functions working again
# Rewrite the per-node time-cost values in rows 8, 10 and 12 of the
# "All_Short_Solutions" sheet (columns F.. onward hold the per-node values;
# rows got shorter because there are now fewer nodes being costed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 --------------------------------------------------------------
$row8 = @{
    "F8" = 44.0;  "G8" = 49.0;  "H8" = 57.0;  "I8" = 50.0;  "J8" = 36.0
    "K8" = 31.0;  "L8" = 30.0;  "M8" = 28.0;  "N8" = 29.0;  "O8" = 27.0
    "P8" = 26.0;  "Q8" = 19.0;  "R8" = 18.0;  "S8" = 11.0;  "T8" = 33.0
    "U8" = 34.0;  "V8" = 4.0;   "W8" = 5.0;   "X8" = 40.0;  "Y8" = 39.0
    "Z8" = 6.0;   "AA8" = 7.0;  "AB8" = -1.0
}
foreach ($addr in $row8.Keys) {
    $ws.Range($addr).Value = $row8[$addr]
}
$row8Clear = @("AC8","AD8","AE8","AF8","AG8","AH8","AI8","AJ8","AK8","AL8","AM8","AN8","AO8","AP8","AQ8","AR8","AS8","AT8")
foreach ($addr in $row8Clear) {
    $ws.Range($addr).ClearContents()
}

# --- Row 10 ---------------------------------------------------------------
$row10 = @{
    "P10" = 50.0; "Q10" = 57.0; "R10" = 49.0; "S10" = 44.0
    "T10" = 39.0; "U10" = 40.0; "V10" = -1.0
}
foreach ($addr in $row10.Keys) {
    $ws.Range($addr).Value = $row10[$addr]
}
$row10Clear = @("W10","X10","Y10","Z10","AA10","AB10")
foreach ($addr in $row10Clear) {
    $ws.Range($addr).ClearContents()
}

# --- Row 12 ---------------------------------------------------------------
$row12 = @{
    "F12" = 36.0; "G12" = 34.0; "H12" = 31.0; "I12" = 18.0; "J12" = 26.0
    "K12" = 27.0; "L12" = 29.0; "M12" = 28.0; "N12" = 30.0; "O12" = 50.0
    "P12" = 57.0; "Q12" = 49.0; "R12" = 44.0; "S12" = 40.0; "T12" = -1.0
}
foreach ($addr in $row12.Keys) {
    $ws.Range($addr).Value = $row12[$addr]
}
$row12Clear = @("U12","V12","W12","X12","Y12","Z12")
foreach ($addr in $row12Clear) {
    $ws.Range($addr).ClearContents()
}
